# Ran code for averaged intensities on spiral schemes.
# Adds three new rotation schemes ("Spiral-90deg-10rot-5space",
# "Spiral-90deg-15rot-5space" and "Spiral-90deg-10rot-3space") to the
# averaged-intensities table, and re-populates the existing rows with the
# newly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows to make room for the 3 new Spiral-scheme data rows
$ws.Rows.Item(11).Resize(3).Insert()

# Fix up column-A style on the newly inserted rows to match the rest of the index column
$ws.Range("A10").Copy()
$ws.Range("A11:A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9975163270514803
$ws.Range("D10").Value = 0.9935294117647059
$ws.Range("E10").Value = 0.9976763069520145
$ws.Range("F10").Value = 0.9975163270514803
$ws.Range("G10").Value = 0.9964705882352941
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0.9976470588235294
$ws.Range("J10").Value = 0.9935294117647059
$ws.Range("K10").Value = 0.9956028593583601
$ws.Range("L10").Value = 0.9965595932049203
$ws.Range("M10").Value = 0.9971399488045041

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9990747111408651
$ws.Range("D11").Value = 0.9991682736143837
$ws.Range("E11").Value = 0.9959366465254276
$ws.Range("F11").Value = 0.9990747111408651
$ws.Range("G11").Value = 0.9994270277724731
$ws.Range("H11").Value = 0.9965092648722953
$ws.Range("I11").Value = 0.9959301583071776
$ws.Range("J11").Value = 0.9991682736143837
$ws.Range("K11").Value = 0.9975524600699056
$ws.Range("L11").Value = 0.9983135856053854
$ws.Range("M11").Value = 0.9976743470387704

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9990965650535065
$ws.Range("D12").Value = 0.9992010866660402
$ws.Range("E12").Value = 0.9959265624668212
$ws.Range("F12").Value = 0.9990965650535065
$ws.Range("G12").Value = 0.9994228328970023
$ws.Range("H12").Value = 0.9965034951466322
$ws.Range("I12").Value = 0.9959221781012233
$ws.Range("J12").Value = 0.9992010866660402
$ws.Range("K12").Value = 0.9975638245664307
$ws.Range("L12").Value = 0.9983301948099687
$ws.Range("M12").Value = 0.9976787867218709

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9990769107758228
$ws.Range("D13").Value = 0.9992095283665813
$ws.Range("E13").Value = 0.9959348334331688
$ws.Range("F13").Value = 0.9990769107758228
$ws.Range("G13").Value = 0.999423094103247
$ws.Range("H13").Value = 0.9965114376054066
$ws.Range("I13").Value = 0.9959308494181962
$ws.Range("J13").Value = 0.9992095283665813
$ws.Range("K13").Value = 0.997572180899875
$ws.Range("L13").Value = 0.9983245458378489
$ws.Range("M13").Value = 0.9976811089504037

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0.9835320000000002
$ws.Range("E14").Value = 0.9979199999999997
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.9995839999999997
$ws.Range("H14").Value = 0.9981239999999996
$ws.Range("I14").Value = 0.9987520000000008
$ws.Range("J14").Value = 0.9835320000000002
$ws.Range("K14").Value = 0.990726
$ws.Range("L14").Value = 0.995363
$ws.Range("M14").Value = 0.9963186666666667

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0.97
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 0.97
$ws.Range("K15").Value = 0.985
$ws.Range("L15").Value = 0.9924999999999999
$ws.Range("M15").Value = 0.995

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9993286389760015
$ws.Range("D16").Value = 0.9831596593152018
$ws.Range("E16").Value = 0.998872056832001
$ws.Range("F16").Value = 0.9993286389760015
$ws.Range("G16").Value = 0.9986210555904008
$ws.Range("H16").Value = 0.9993782108159984
$ws.Range("I16").Value = 0.9977968009215982
$ws.Range("J16").Value = 0.9831596593152018
$ws.Range("K16").Value = 0.9910158580736014
$ws.Range("L16").Value = 0.9951722485248015
$ws.Range("M16").Value = 0.9961927370752003

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9984636164895592
$ws.Range("D17").Value = 0.9990816090398335
$ws.Range("E17").Value = 0.9964945944622892
$ws.Range("F17").Value = 0.9984636164895592
$ws.Range("G17").Value = 0.9971011832630307
$ws.Range("H17").Value = 0.9986177966212647
$ws.Range("I17").Value = 0.9942589405080793
$ws.Range("J17").Value = 0.9990816090398335
$ws.Range("K17").Value = 0.9977881017510613
$ws.Range("L17").Value = 0.9981258591203102
$ws.Range("M17").Value = 0.9973362900640094

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9977182519561624
$ws.Range("D18").Value = 0.9983720930232558
$ws.Range("E18").Value = 0.9960889352813616
$ws.Range("F18").Value = 0.9977182519561624
$ws.Range("G18").Value = 0.9967262564741518
$ws.Range("H18").Value = 0.9986110584363935
$ws.Range("I18").Value = 0.9936869216454325
$ws.Range("J18").Value = 0.9983720930232558
$ws.Range("K18").Value = 0.9972305141523087
$ws.Range("L18").Value = 0.9974743830542355
$ws.Range("M18").Value = 0.9968672528027929

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9979522955577034
$ws.Range("D19").Value = 0.9993523139911924
$ws.Range("E19").Value = 0.9970507227452089
$ws.Range("F19").Value = 0.9979522955577034
$ws.Range("G19").Value = 0.9983300879800956
$ws.Range("H19").Value = 0.9979914168531766
$ws.Range("I19").Value = 0.9938847166986744
$ws.Range("J19").Value = 0.9993523139911924
$ws.Range("K19").Value = 0.9982015183682007
$ws.Range("L19").Value = 0.998076906962952
$ws.Range("M19").Value = 0.9974269256376752

